$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.03562113421745039
$ws.Range("D2").Value = 0.156229611855025
$ws.Range("E2").Value = 0.1225234119618293
$ws.Range("F2").Value = 3.164852129049535
$ws.Range("G2").Value = 0.002400743933228427
$ws.Range("I2").Value = 2.470448866583922
$ws.Range("M2").Value = 2.106466521024117
$ws.Range("N2").Value = 1.35901728559125

$ws.Range("C3").Value = 0.03273057732545226
$ws.Range("D3").Value = 0.1426978178676563
$ws.Range("E3").Value = 0.1074248718158941
$ws.Range("F3").Value = 2.913608175433694
$ws.Range("G3").Value = 0.002413370137977345
$ws.Range("I3").Value = 2.27303571660488
$ws.Range("M3").Value = 1.868222204624544
$ws.Range("N3").Value = 1.273358158093657

$ws.Range("C4").Value = 0.03100276902115695
$ws.Range("D4").Value = 0.1345767623620162
$ws.Range("E4").Value = 0.09811675115843599
$ws.Range("F4").Value = 2.762657873889196
$ws.Range("G4").Value = 0.002421489262605471
$ws.Range("I4").Value = 2.154420543968172
$ws.Range("M4").Value = 1.722656868044112
$ws.Range("N4").Value = 1.221123521883698

$ws.Range("C5").Value = 0.03030996141485787
$ws.Range("D5").Value = 0.131312308262693
$ws.Range("E5").Value = 0.0943143315991648
$ws.Range("F5").Value = 2.701943552075733
$ws.Range("G5").Value = 0.002424890636847941
$ws.Range("I5").Value = 2.106709905902335
$ws.Range("M5").Value = 1.663505478158896
$ws.Range("N5").Value = 1.199930509835099

$ws.Range("C6").Value = 0.03019558775420705
$ws.Range("D6").Value = 0.1307728990121007
$ws.Range("E6").Value = 0.09368238676326968
$ws.Range("F6").Value = 2.691909267140801
$ws.Range("G6").Value = 0.002425461052072966
$ws.Range("I6").Value = 2.098824622037171
$ws.Range("M6").Value = 1.653693196820456
$ws.Range("N6").Value = 1.196417116903092

$ws.Range("C7").Value = 0.03099338056129852
$ws.Range("D7").Value = 0.1345325577030678
$ws.Range("E7").Value = 0.09806550777720702
$ws.Range("F7").Value = 2.761835868219691
$ws.Range("G7").Value = 0.002421534758538163
$ws.Range("I7").Value = 2.153774602206312
$ws.Range("M7").Value = 1.721858469403372
$ws.Range("N7").Value = 1.220837325613218

$ws.Range("C8").Value = 0.03461438913126358
$ws.Range("D8").Value = 0.151523575965399
$ws.Range("E8").Value = 0.1173253187944709
$ws.Range("F8").Value = 3.077514274894781
$ws.Range("G8").Value = 0.002405021719130892
$ws.Range("I8").Value = 2.401825141204142
$ws.Range("M8").Value = 2.024163024651159
$ws.Range("N8").Value = 1.329408598983633

$ws.Range("C9").Value = 0.0421130934652183
$ws.Range("D9").Value = 0.1864348391532644
$ws.Range("E9").Value = 0.1547909305856052
$ws.Range("F9").Value = 3.72447873785535
$ws.Range("G9").Value = 0.002375521359838912
$ws.Range("I9").Value = 2.910140354498338
$ws.Range("M9").Value = 2.623302999363773
$ws.Range("N9").Value = 1.545086185056249

$ws.Range("C10").Value = 0.04790103716578642
$ws.Range("D10").Value = 0.2132053107636693
$ws.Range("E10").Value = 0.1821290808529028
$ws.Range("F10").Value = 4.219182146296873
$ws.Range("G10").Value = 0.002355565786795926
$ws.Range("I10").Value = 3.298805766249416
$ws.Range("M10").Value = 3.06829035325751
$ws.Range("N10").Value = 1.705134347064757

$ws.Range("C11").Value = 0.05060283496126772
$ws.Range("D11").Value = 0.225661665631776
$ws.Range("E11").Value = 0.1945248136804167
$ws.Range("F11").Value = 4.448981203440098
$ws.Range("G11").Value = 0.002346852125375361
$ws.Range("I11").Value = 3.479346337402632
$ws.Range("M11").Value = 3.271987933714712
$ws.Range("N11").Value = 1.778271474700631

$ws.Range("C12").Value = 0.05163652835884136
$ws.Range("D12").Value = 0.2304214882888402
$ws.Range("E12").Value = 0.1992128554298347
$ws.Range("F12").Value = 4.536729304073049
$ws.Range("G12").Value = 0.002343604177584956
$ws.Range("I12").Value = 3.548285162186403
$ws.Range("M12").Value = 3.349323671978766
$ws.Range("N12").Value = 1.806012349716838

$ws.Range("C13").Value = 0.05141342222549383
$ws.Range("D13").Value = 0.2293944220505182
$ws.Range("E13").Value = 0.1982034683846621
$ws.Range("F13").Value = 4.517798070352057
$ws.Range("G13").Value = 0.002344301391248495
$ws.Range("I13").Value = 3.533411939704308
$ws.Range("M13").Value = 3.332658864706303
$ws.Range("N13").Value = 1.800035867523434

$ws.Range("C14").Value = 0.0506876614671512
$ws.Range("D14").Value = 0.2260523831815817
$ws.Range("E14").Value = 0.1949106216756746
$ws.Range("F14").Value = 4.45618543862679
$ws.Range("G14").Value = 0.002346583881934587
$ws.Range("I14").Value = 3.485006306337937
$ws.Range("M14").Value = 3.278346284089253
$ws.Range("N14").Value = 1.780552834267695

$ws.Range("C15").Value = 0.05024451076228331
$ws.Range("D15").Value = 0.2240109574370877
$ws.Range("E15").Value = 0.192892879210369
$ws.Range("F15").Value = 4.418542098182115
$ws.Range("G15").Value = 0.002347988690150536
$ws.Range("I15").Value = 3.455432015256946
$ws.Range("M15").Value = 3.24510481509887
$ws.Range("N15").Value = 1.768624770174711

$ws.Range("C16").Value = 0.04772591121520975
$ws.Range("D16").Value = 0.2123971003196061
$ws.Range("E16").Value = 0.1813181640137671
$ws.Range("F16").Value = 4.204263685414219
$ws.Range("G16").Value = 0.002356142519152172
$ws.Range("I16").Value = 3.287085126205966
$ws.Range("M16").Value = 3.05500525445197
$ws.Range("N16").Value = 1.700361132632906

$ws.Range("C17").Value = 0.04619893952965981
$ws.Range("D17").Value = 0.2053456539420324
$ws.Range("E17").Value = 0.1742069775508739
$ws.Range("F17").Value = 4.074059899465851
$ws.Range("G17").Value = 0.00236123745969314
$ws.Range("I17").Value = 3.184790769895443
$ws.Range("M17").Value = 2.938722724572358
$ws.Range("N17").Value = 1.658566720375262

$ws.Range("C18").Value = 0.04532710228045289
$ws.Range("D18").Value = 0.2013158695717721
$ws.Range("E18").Value = 0.1701130031250173
$ws.Range("F18").Value = 3.999615029725817
$ws.Range("G18").Value = 0.002364202256895965
$ws.Range("I18").Value = 3.126303117286938
$ws.Range("M18").Value = 2.871958607549288
$ws.Range("N18").Value = 1.634558898975939

$ws.Range("C19").Value = 0.04503299825324802
$ws.Range("D19").Value = 0.199955838020486
$ws.Range("E19").Value = 0.1687262037190322
$ws.Range("F19").Value = 3.974484367971684
$ws.Range("G19").Value = 0.00236521200103993
$ws.Range("I19").Value = 3.106559161863316
$ws.Range("M19").Value = 2.849373270436729
$ws.Range("N19").Value = 1.626435689950654

$ws.Range("C20").Value = 0.04636081725713836
$ws.Range("D20").Value = 0.2060935790925669
$ws.Range("E20").Value = 0.1749643705573263
$ws.Range("F20").Value = 4.087873953701148
$ws.Range("G20").Value = 0.002360691546888474
$ws.Range("I20").Value = 3.195643772806051
$ws.Range("M20").Value = 2.951088810370749
$ws.Range("N20").Value = 1.663012591208002

$ws.Range("C21").Value = 0.05090054224773155
$ws.Range("D21").Value = 0.2270328342441701
$ws.Range("E21").Value = 0.1958779730786375
$ws.Range("F21").Value = 4.474262437681091
$ws.Range("G21").Value = 0.002345912060441745
$ws.Range("I21").Value = 3.49920840382498
$ws.Range("M21").Value = 3.294293639757768
$ws.Range("N21").Value = 1.786274256539457

$ws.Range("C22").Value = 0.05392952670482032
$ws.Range("D22").Value = 0.2409691384552275
$ws.Range("E22").Value = 0.2095114726814558
$ws.Range("F22").Value = 4.731055899553553
$ws.Range("G22").Value = 0.002336554003770086
$ws.Range("I22").Value = 3.700956906678385
$ws.Range("M22").Value = 3.519771492498137
$ws.Range("N22").Value = 1.867096790430736

$ws.Range("C23").Value = 0.05230699389714744
$ws.Range("D23").Value = 0.2335071070576191
$ws.Range("E23").Value = 0.2022382347743843
$ws.Range("F23").Value = 4.593594938336764
$ws.Range("G23").Value = 0.00234152123362262
$ws.Range("I23").Value = 3.592961349025956
$ws.Range("M23").Value = 3.399316436483502
$ws.Range("N23").Value = 1.823936799839828

$ws.Range("C24").Value = 0.04628761357761846
$ws.Range("D24").Value = 0.2057553670117045
$ws.Range("E24").Value = 0.1746219707946324
$ws.Range("F24").Value = 4.081627342870831
$ws.Range("G24").Value = 0.002360938242903164
$ws.Range("I24").Value = 3.190736127570943
$ws.Range("M24").Value = 2.945497829364768
$ws.Range("N24").Value = 1.661002549481623

$ws.Range("C25").Value = 0.04003824695784886
$ws.Range("D25").Value = 0.1768047496176735
$ws.Range("E25").Value = 0.1446882789437396
$ws.Range("F25").Value = 3.546230227342079
$ws.Range("G25").Value = 0.002383197416242242
$ws.Range("I25").Value = 2.770095942708423
$ws.Range("M25").Value = 2.460447760820898
$ws.Range("N25").Value = 1.486453699468257
